$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I1 ("I0") and J1 ("IF"), matching the style (s="1") of H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-15: fill columns I and J with the corresponding values
$values = @(
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(7, 8),
    @(9, 9),
    @(6, 8),
    @(6, 7),
    @(9, 9),
    @(8, 8),
    @(6, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
